# Scheduled-runner refresh: pulls latest Universalis market-board
# averages and re-derives the LevePrice / LeveProfit columns (H:N)
# for every Disciple of the Hand job sheet in this workbook.

$wb = $excel.ActiveWorkbook

# --- ALC sheet: refreshed price/profit figures ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 415.2857
$ws.Range("J2").Value = 749
$ws.Range("L2").Value = 749
$ws.Range("N2").Value = -975
$ws.Range("H9").Value = 216.69565
$ws.Range("I9").Value = 206.8
$ws.Range("J9").Value = 235.25
$ws.Range("K9").Value = 206.8
$ws.Range("L9").Value = 235.25
$ws.Range("M9").Value = -37.80000000000001
$ws.Range("N9").Value = -573.25
$ws.Range("H121").Value = 1167.4375
$ws.Range("I121").Value = 1333.3334
$ws.Range("J121").Value = 1129.1538
$ws.Range("K121").Value = 4000.0002
$ws.Range("L121").Value = 3387.4614
$ws.Range("M121").Value = -2253.0002
$ws.Range("N121").Value = -6881.4614
$ws.Range("H137").Value = 2496.0278
$ws.Range("I137").Value = 1358.3654
$ws.Range("J137").Value = 5453.95
$ws.Range("K137").Value = 4075.0962
$ws.Range("L137").Value = 16361.85
$ws.Range("M137").Value = -1525.0962
$ws.Range("N137").Value = -21461.85
$ws.Range("H138").Value = 3480.3196
$ws.Range("I138").Value = 1751.4333
$ws.Range("J138").Value = 4254.4478
$ws.Range("K138").Value = 5254.2999
$ws.Range("L138").Value = 12763.3434
$ws.Range("M138").Value = -114.2999
$ws.Range("N138").Value = -23043.3434
$ws.Range("H141").Value = 2638.423
$ws.Range("I141").Value = 2067.1667
$ws.Range("J141").Value = 3923.75
$ws.Range("K141").Value = 6201.500100000001
$ws.Range("L141").Value = 11771.25
$ws.Range("M141").Value = -1021.500100000001
$ws.Range("N141").Value = -22131.25

# --- ARM sheet: refreshed price/profit figures ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 5260
$ws.Range("I25").Value = 5260
$ws.Range("K25").Value = 5260
$ws.Range("M25").Value = -4858
$ws.Range("H32").Value = 5558.737
$ws.Range("I32").Value = 5298.7236
$ws.Range("K32").Value = 5298.7236
$ws.Range("M32").Value = -5011.7236
$ws.Range("H61").Value = 6527.222
$ws.Range("I61").Value = 3609.0645
$ws.Range("J61").Value = 12988.857
$ws.Range("K61").Value = 3609.0645
$ws.Range("L61").Value = 12988.857
$ws.Range("M61").Value = -3397.0645
$ws.Range("N61").Value = -13412.857
$ws.Range("H122").Value = 25002244
$ws.Range("I122").Value = 2222
$ws.Range("K122").Value = 6666
$ws.Range("M122").Value = -4216
$ws.Range("H132").Value = 4953.114
$ws.Range("I132").Value = 1578.4584
$ws.Range("J132").Value = 9002.700000000001
$ws.Range("K132").Value = 4735.3752
$ws.Range("L132").Value = 27008.1
$ws.Range("M132").Value = -2205.3752
$ws.Range("N132").Value = -32068.1
$ws.Range("H136").Value = 6527.222
$ws.Range("I136").Value = 3609.0645
$ws.Range("J136").Value = 12988.857
$ws.Range("K136").Value = 10827.1935
$ws.Range("L136").Value = 38966.571
$ws.Range("M136").Value = -8277.193499999999
$ws.Range("N136").Value = -44066.571
$ws.Range("H139").Value = 49143.332
$ws.Range("J139").Value = 49143.332
$ws.Range("L139").Value = 49143.332
$ws.Range("N139").Value = -59423.332

# --- BSM sheet: refreshed price/profit figures ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1153.3
$ws.Range("I20").Value = 1154
$ws.Range("J20").Value = 1152.25
$ws.Range("K20").Value = 1154
$ws.Range("L20").Value = 1152.25
$ws.Range("M20").Value = -907
$ws.Range("N20").Value = -1646.25
$ws.Range("H81").Value = 16328.385
$ws.Range("J81").Value = 16328.385
$ws.Range("L81").Value = 16328.385
$ws.Range("N81").Value = -18450.385
$ws.Range("H84").Value = 16328.385
$ws.Range("J84").Value = 16328.385
$ws.Range("L84").Value = 48985.155
$ws.Range("N84").Value = -59593.155
$ws.Range("H134").Value = 3715.1904
$ws.Range("I134").Value = 3588.5293
$ws.Range("K134").Value = 10765.5879
$ws.Range("M134").Value = -8230.5879

# --- CRP sheet: refreshed price/profit figures ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2211.2158
$ws.Range("I31").Value = 1487.3151
$ws.Range("J31").Value = 5734.2
$ws.Range("K31").Value = 1487.3151
$ws.Range("L31").Value = 5734.2
$ws.Range("M31").Value = -1192.3151
$ws.Range("N31").Value = -6324.2
$ws.Range("H34").Value = 2211.2158
$ws.Range("I34").Value = 1487.3151
$ws.Range("J34").Value = 5734.2
$ws.Range("K34").Value = 1487.3151
$ws.Range("L34").Value = 5734.2
$ws.Range("M34").Value = -1285.3151
$ws.Range("N34").Value = -6138.2
$ws.Range("H39").Value = 6707.1
$ws.Range("I39").Value = 3585.6667
$ws.Range("K39").Value = 3585.6667
$ws.Range("M39").Value = -3194.6667
$ws.Range("H49").Value = 6707.1
$ws.Range("I49").Value = 3585.6667
$ws.Range("K49").Value = 3585.6667
$ws.Range("M49").Value = -3403.6667
$ws.Range("H99").Value = 2254.0833
$ws.Range("I99").Value = 2215
$ws.Range("J99").Value = 2449.5
$ws.Range("K99").Value = 2215
$ws.Range("L99").Value = 2449.5
$ws.Range("M99").Value = -717
$ws.Range("N99").Value = -5445.5
$ws.Range("H105").Value = 572.5714
$ws.Range("I105").Value = 572.5714
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 572.5714
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 1174.4286
$ws.Range("N105").ClearContents()
$ws.Range("H126").Value = 2254.0833
$ws.Range("I126").Value = 2215
$ws.Range("J126").Value = 2449.5
$ws.Range("K126").Value = 6645
$ws.Range("L126").Value = 7348.5
$ws.Range("M126").Value = -4175
$ws.Range("N126").Value = -12288.5

# --- CUL sheet: refreshed price/profit figures ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 295
$ws.Range("I17").Value = 280
$ws.Range("K17").Value = 840
$ws.Range("M17").Value = -671
$ws.Range("H39").Value = 13212.75
$ws.Range("J39").Value = 13212.75
$ws.Range("L39").Value = 39638.25
$ws.Range("N39").Value = -40226.25
$ws.Range("H55").Value = 6333.3335
$ws.Range("J55").Value = 6333.3335
$ws.Range("L55").Value = 19000.0005
$ws.Range("N55").Value = -19354.0005
$ws.Range("H127").Value = 3694.9395
$ws.Range("J127").Value = 3694.9395
$ws.Range("L127").Value = 11084.8185
$ws.Range("N127").Value = -21004.8185
$ws.Range("H131").Value = 529.3200000000001
$ws.Range("I131").Value = 288.1579
$ws.Range("K131").Value = 864.4737
$ws.Range("M131").Value = 4175.5263
$ws.Range("H132").Value = 1867.3125
$ws.Range("I132").Value = 1999.6666
$ws.Range("J132").Value = 1697.1428
$ws.Range("K132").Value = 17996.9994
$ws.Range("L132").Value = 15274.2852
$ws.Range("M132").Value = -15466.9994
$ws.Range("N132").Value = -20334.2852

# --- GSM sheet: refreshed price/profit figures ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 3087.25
$ws.Range("I36").Value = 2279.5
$ws.Range("J36").Value = 3895
$ws.Range("K36").Value = 2279.5
$ws.Range("L36").Value = 3895
$ws.Range("M36").Value = -1794.5
$ws.Range("N36").Value = -4865
$ws.Range("H102").Value = 3506.4666
$ws.Range("I102").Value = 3219.2188
$ws.Range("J102").Value = 4213.5386
$ws.Range("K102").Value = 3219.2188
$ws.Range("L102").Value = 4213.5386
$ws.Range("M102").Value = -1597.2188
$ws.Range("N102").Value = -7457.5386
$ws.Range("H122").Value = 4195.4287
$ws.Range("I122").Value = 5946.5454
$ws.Range("J122").Value = 2269.2
$ws.Range("K122").Value = 17839.6362
$ws.Range("L122").Value = 6807.599999999999
$ws.Range("M122").Value = -15389.6362
$ws.Range("N122").Value = -11707.6
$ws.Range("H123").Value = 27300
$ws.Range("J123").Value = 27300
$ws.Range("L123").Value = 27300
$ws.Range("N123").Value = -32200
$ws.Range("H135").Value = 63750
$ws.Range("J135").Value = 63750
$ws.Range("L135").Value = 63750
$ws.Range("N135").Value = -73890

# --- LTW sheet: refreshed price/profit figures ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4154.2144
$ws.Range("I7").Value = 3935.5
$ws.Range("J7").Value = 4445.8335
$ws.Range("K7").Value = 3935.5
$ws.Range("L7").Value = 4445.8335
$ws.Range("M7").Value = -3823.5
$ws.Range("N7").Value = -4669.8335
$ws.Range("H22").Value = 954.7143
$ws.Range("J22").Value = 756.6
$ws.Range("L22").Value = 756.6
$ws.Range("N22").Value = -1346.6
$ws.Range("H27").Value = 954.7143
$ws.Range("J27").Value = 756.6
$ws.Range("L27").Value = 756.6
$ws.Range("N27").Value = -970.6
$ws.Range("H38").Value = 19990.2
$ws.Range("J38").Value = 19990.2
$ws.Range("L38").Value = 19990.2
$ws.Range("N38").Value = -20810.2
$ws.Range("H40").Value = 3992.125
$ws.Range("I40").Value = 3724.32
$ws.Range("K40").Value = 3724.32
$ws.Range("M40").Value = -3588.32
$ws.Range("H61").Value = 51068
$ws.Range("I61").Value = 60981.6
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 60981.6
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -60779.6
$ws.Range("N61").Value = -1904
$ws.Range("H113").Value = 51068
$ws.Range("I113").Value = 60981.6
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 60981.6
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -58811.6
$ws.Range("N113").Value = -5840
$ws.Range("H126").Value = 4154.2144
$ws.Range("I126").Value = 3935.5
$ws.Range("J126").Value = 4445.8335
$ws.Range("K126").Value = 11806.5
$ws.Range("L126").Value = 13337.5005
$ws.Range("M126").Value = -9336.5
$ws.Range("N126").Value = -18277.5005

# --- WVR sheet: refreshed price/profit figures ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5109.826
$ws.Range("I122").Value = 2936.5
$ws.Range("J122").Value = 8490.556
$ws.Range("K122").Value = 8809.5
$ws.Range("L122").Value = 25471.668
$ws.Range("M122").Value = -6359.5
$ws.Range("N122").Value = -30371.668

